# Adding referral profile and modifying current profiles
#
# - Bump the "Date" value on the Metadata sheet.
# - Duplicate the "Include #0" worksheet to create a new "Include #1"
#   worksheet (this preserves the existing cell styles/column widths/blank
#   cells exactly, which a plain value-by-value copy would not do).
# - "Include #0" keeps its System URI row but the value is changed to the
#   HL7 ActCode system.
# - The new "Include #1" sheet keeps the System URI value that
#   "Include #0" used to have (the Eswatini encounter-classification
#   CodeSystem).

$wb = $excel.ActiveWorkbook

# --- Metadata sheet: refresh generation Date ---------------------------
$metaSheet = $wb.Worksheets.Item(1)
$metaSheet.Range("B8").Value = "2025-08-01T12:25:19+00:00"

# --- Duplicate "Include #0" as the new "Include #1" sheet --------------
$include0 = $wb.Worksheets.Item(2)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$include0.Copy([System.Reflection.Missing]::Value, $lastSheet)

$include1 = $wb.Worksheets.Item($wb.Worksheets.Count)
$include1.Name = "Include #1"

# --- Modify the current "Include #0" profile's System URI --------------
$include0.Range("B4").Value = "http://terminology.hl7.org/CodeSystem/v3-ActCode"

# --- Restore original active sheet/selection ----------------------------
$metaSheet.Activate()
